$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "STORJ"
$ws.Range("B29").Value = 1649
$ws.Range("C29").Value = 0.658
$ws.Range("D29").Value = "Storage"

$ws.Range("A30").Select()
